$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3363324105739594
$ws.Range("B1").Value = 0.2017912417650223
$ws.Range("C1").Value = 0.2689874768257141
$ws.Range("D1").Value = 3.49756908416748
$ws.Range("E1").Value = 3.939542293548584
